$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Insert a new row above row 6 (indicators were renumbered/corrected; a new
# row for INDICATOR_11 is inserted and everything below shifts down by one).
$ws.Rows.Item(6).Insert()

# Populate the new row 6 following the same pattern as its sibling rows.
$ws.Range("A6").Value = "CREATE/MODIFY"
$ws.Range("B6").Value = "LIB_EWS_BE"
$ws.Range("C6").Value = "INDICATOR_11"
$ws.Range("E6").Value = "String"

# Restore the view state: scrolled near the top, with C6 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C6").Select()
